# LipnoV2 TOPOLOGIE DONE 24.2.15:21
# Fill in the missing start/stop time for the last tracked interval (row 255)
# and tag it with the "TOPOLOGIE DONE" marker, matching the day's other rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start (C) / stop (D) times for row 255 -> 14:49 and 15:19
$ws.Range("C255").Value = 0.61736111111111114
$ws.Range("D255").Value = 0.6381944444444444

# Mark this entry as the point "TOPOLOGIE" work was finished
$ws.Range("J255").Value = "TOPOLOGIE DONE"

# Leave the cursor where the author ended up while reviewing the change
$ws.Range("L259").Select()
